$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether to force Text number format
# (needed for numeric-looking strings like "1.010" / "26.363.58" so Excel
# keeps them as literal text instead of re-parsing them as numbers).
$updates = @(
    @{ Cell = 'D2'; Value = '26.363.58'; Text = $true },
    @{ Cell = 'E2'; Value = '  +0.42%  '; Text = $true },
    @{ Cell = 'D3'; Value = '1.689.65'; Text = $true },
    @{ Cell = 'E3'; Value = '  +0.01%  '; Text = $true },
    @{ Cell = 'D4'; Value = '1.010'; Text = $true },
    @{ Cell = 'E4'; Value = '  +0.40%  '; Text = $true },
    @{ Cell = 'D5'; Value = '218.81'; Text = $true },
    @{ Cell = 'E5'; Value = '  -0.19%  '; Text = $true },
    @{ Cell = 'D6'; Value = '0.5459'; Text = $true },
    @{ Cell = 'E6'; Value = '  +4.06%  '; Text = $true },
    @{ Cell = 'E7'; Value = '  +0.33%  '; Text = $true },
    @{ Cell = 'D8'; Value = '0.2745'; Text = $true },
    @{ Cell = 'E8'; Value = '  +1.85%  '; Text = $true },
    @{ Cell = 'D9'; Value = '0.06457'; Text = $true },
    @{ Cell = 'E9'; Value = '  +0.24%  '; Text = $true },
    @{ Cell = 'D10'; Value = '21.99'; Text = $true },
    @{ Cell = 'E10'; Value = '  -0.42%  '; Text = $true },
    @{ Cell = 'D11'; Value = '0.07689'; Text = $true },
    @{ Cell = 'E11'; Value = '  +3.09%  '; Text = $true },
    @{ Cell = 'D12'; Value = '1.702.13'; Text = $true },
    @{ Cell = 'E12'; Value = '  +0.53%  '; Text = $true },
    @{ Cell = 'D13'; Value = '4.533'; Text = $true },
    @{ Cell = 'E13'; Value = '  -0.47%  '; Text = $true },
    @{ Cell = 'D14'; Value = '0.5832'; Text = $true },
    @{ Cell = 'E14'; Value = '  -0.49%  '; Text = $true },
    @{ Cell = 'D15'; Value = '0.000008383'; Text = $true },
    @{ Cell = 'E15'; Value = '  -1.98%  '; Text = $true },
    @{ Cell = 'D16'; Value = '65.25'; Text = $true },
    @{ Cell = 'E16'; Value = '  +0.90%  '; Text = $true },
    @{ Cell = 'D17'; Value = '26.397.65'; Text = $true },
    @{ Cell = 'E17'; Value = '  +0.30%  '; Text = $true },
    @{ Cell = 'D18'; Value = '4.941'; Text = $true },
    @{ Cell = 'E18'; Value = '  -0.58%  '; Text = $true },
    @{ Cell = 'E19'; Value = '  +0.40%  '; Text = $true },
    @{ Cell = 'D20'; Value = '10.98'; Text = $true },
    @{ Cell = 'E20'; Value = '  +1.14%  '; Text = $true },
    @{ Cell = 'D21'; Value = '191.54'; Text = $true },
    @{ Cell = 'E21'; Value = '  +0.49%  '; Text = $true },
    @{ Cell = 'D22'; Value = '6.257'; Text = $true },
    @{ Cell = 'E22'; Value = '  +0.31%  '; Text = $true },
    @{ Cell = 'E23'; Value = '  +0.38%  '; Text = $true },
    @{ Cell = 'D24'; Value = '149.54'; Text = $true },
    @{ Cell = 'E24'; Value = '  +2.90%  '; Text = $true },
    @{ Cell = 'D25'; Value = '0.1323'; Text = $true },
    @{ Cell = 'E25'; Value = '  +6.40%  '; Text = $true },
    @{ Cell = 'D26'; Value = '7.890'; Text = $true },
    @{ Cell = 'E26'; Value = '  +2.71%  '; Text = $true },
    @{ Cell = 'E27'; Value = '  -0.88%  '; Text = $true },
    @{ Cell = 'D28'; Value = '0.06356'; Text = $true },
    @{ Cell = 'E28'; Value = '  -4.61%  '; Text = $true },
    @{ Cell = 'E29'; Value = '  +3.82%  '; Text = $true },
    @{ Cell = 'D30'; Value = '1.328'; Text = $true },
    @{ Cell = 'E30'; Value = '  -0.25%  '; Text = $true },
    @{ Cell = 'D31'; Value = '3.598'; Text = $true },
    @{ Cell = 'E31'; Value = '  -0.16%  '; Text = $true },
    @{ Cell = 'D32'; Value = '3.587'; Text = $true },
    @{ Cell = 'E32'; Value = '  +0.89%  '; Text = $true },
    @{ Cell = 'D33'; Value = '1.687'; Text = $true },
    @{ Cell = 'E33'; Value = '  +1.22%  '; Text = $true },
    @{ Cell = 'D34'; Value = '1.044'; Text = $true },
    @{ Cell = 'E34'; Value = '  +1.52%  '; Text = $true },
    @{ Cell = 'D35'; Value = '0.6151'; Text = $true },
    @{ Cell = 'E35'; Value = '  -0.98%  '; Text = $true },
    @{ Cell = 'E36'; Value = '  +0.98%  '; Text = $true },
    @{ Cell = 'D37'; Value = '2.709'; Text = $true },
    @{ Cell = 'E37'; Value = '  +0.12%  '; Text = $true },
    @{ Cell = 'D38'; Value = '6.272'; Text = $true },
    @{ Cell = 'E38'; Value = '  -0.01%  '; Text = $true },
    @{ Cell = 'D39'; Value = '1.118.47'; Text = $true },
    @{ Cell = 'E39'; Value = '  +1.40%  '; Text = $true },
    @{ Cell = 'D40'; Value = '0.01632'; Text = $true },
    @{ Cell = 'E40'; Value = '  +0.78%  '; Text = $true },
    @{ Cell = 'D41'; Value = '0.8789'; Text = $true },
    @{ Cell = 'E41'; Value = '  -0.23%  '; Text = $true },
    @{ Cell = 'E42'; Value = '  -0.09%  '; Text = $true },
    @{ Cell = 'D43'; Value = '101.76'; Text = $true },
    @{ Cell = 'E43'; Value = '  +0.87%  '; Text = $true },
    @{ Cell = 'D44'; Value = '1.838.79'; Text = $true },
    @{ Cell = 'E44'; Value = '  +0.08%  '; Text = $true },
    @{ Cell = 'D45'; Value = '57.51'; Text = $true },
    @{ Cell = 'E45'; Value = '  +1.09%  '; Text = $true },
    @{ Cell = 'D46'; Value = '0.00000000108'; Text = $true },
    @{ Cell = 'E46'; Value = '  -2.73%  '; Text = $true },
    @{ Cell = 'E47'; Value = '  +0.75%  '; Text = $true },
    @{ Cell = 'D48'; Value = '8.212'; Text = $true },
    @{ Cell = 'E48'; Value = '  +0.63%  '; Text = $true },
    @{ Cell = 'D49'; Value = '0.05273'; Text = $true },
    @{ Cell = 'E49'; Value = '  +0.20%  '; Text = $true },
    @{ Cell = 'B50'; Value = 'Aptos'; Text = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Text = $false },
    @{ Cell = 'D50'; Value = '6.170'; Text = $true },
    @{ Cell = 'E50'; Value = '  +2.32%  '; Text = $true },
    @{ Cell = 'B51'; Value = 'Mantle'; Text = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; Text = $false },
    @{ Cell = 'D51'; Value = '0.4305'; Text = $true },
    @{ Cell = 'E51'; Value = '  +0.21%  '; Text = $true }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Text) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
